$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two existing tugboats (NB006 / row 7, NB013 / row 14) go into Maintenance.
$ws.Range("B7").Value = "Maintenance"
$ws.Range("B14").Value = "Maintenance"

# Add two new tugboats (NB021, NB022) that are currently Free, each with a
# captain assignment and the same working-hours window as the existing rows.
$ws.Range("A22").Value = "NB021"
$ws.Range("A23").Value = "NB022"

$ws.Range("C22").Value = "CP0021"
$ws.Range("C23").Value = "CP0022"

$ws.Range("B22").Value = "Free"
$ws.Range("B23").Value = "Free"

$ws.Range("D22").Value = 0.375
$ws.Range("E22").Value = 0.75

$ws.Range("D23").Value = 0.375
$ws.Range("E23").Value = 0.75

# Match the plain (unstyled) look of the rest of column B - new cells in a
# column that carries a column-level style otherwise pick that style up.
$plainStyle = $ws.Range("C21").Style()
$ws.Range("B22").Style = $plainStyle
$ws.Range("B23").Style = $plainStyle

# Match the time-column formatting used by the rest of the table.
$ws.Range("D21:E21").Copy()
$ws.Range("D22:E23").PasteSpecial(-4122)

$ws.Rows.Item(22).RowHeight = 15
$ws.Rows.Item(23).RowHeight = 15
